$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H28").Value = 676.5714
$ws.Range("I28").Value = 387.2
$ws.Range("J28").Value = 1400
$ws.Range("K28").Value = 387.2
$ws.Range("L28").Value = 1400
$ws.Range("M28").Value = 97.80000000000001
$ws.Range("N28").Value = -2370
$ws.Range("H43").Value = 857.6
$ws.Range("I43").Value = 895.3333
$ws.Range("J43").Value = 801
$ws.Range("K43").Value = 895.3333
$ws.Range("L43").Value = 801
$ws.Range("M43").Value = -826.3333
$ws.Range("N43").Value = -939
$ws.Range("H87").Value = 30133
$ws.Range("J87").Value = 30133
$ws.Range("L87").Value = 30133
$ws.Range("N87").Value = -32629
$ws.Range("H90").Value = 30133
$ws.Range("J90").Value = 30133
$ws.Range("L90").Value = 90399
$ws.Range("N90").Value = -102879
$ws.Range("H106").Value = 20841078
$ws.Range("I106").Value = 20841078
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 20841078
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -20840447
$ws.Range("N106").ClearContents()
$ws.Range("H136").Value = 57410
$ws.Range("J136").Value = 57410
$ws.Range("L136").Value = 57410
$ws.Range("N136").Value = -67610
$ws.Range("H137").Value = 1411.3871
$ws.Range("I137").Value = 1394
$ws.Range("J137").Value = 1483.8334
$ws.Range("K137").Value = 4182
$ws.Range("L137").Value = 4451.5002
$ws.Range("M137").Value = -1632
$ws.Range("N137").Value = -9551.5002
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H32").Value = 18780.166
$ws.Range("I32").Value = 16707.814
$ws.Range("K32").Value = 16707.814
$ws.Range("M32").Value = -16420.814
$ws.Range("H61").Value = 351645.78
$ws.Range("I61").Value = 10294.214
$ws.Range("J61").Value = 670240.6
$ws.Range("K61").Value = 10294.214
$ws.Range("L61").Value = 670240.6
$ws.Range("M61").Value = -10082.214
$ws.Range("N61").Value = -670664.6
$ws.Range("H136").Value = 351645.78
$ws.Range("I136").Value = 10294.214
$ws.Range("J136").Value = 670240.6
$ws.Range("K136").Value = 30882.642
$ws.Range("L136").Value = 2010721.8
$ws.Range("M136").Value = -28332.642
$ws.Range("N136").Value = -2015821.8
$ws = $wb.Sheets.Item("BSM")
$ws.Range("H80").Value = 326.03845
$ws.Range("I80").Value = 72.85714
$ws.Range("K80").Value = 72.85714
$ws.Range("M80").Value = 925.14286
$ws.Range("H83").Value = 326.03845
$ws.Range("I83").Value = 72.85714
$ws.Range("K83").Value = 364.2857
$ws.Range("M83").Value = 4627.7143
$ws.Range("H86").Value = 1977.7778
$ws.Range("I86").Value = 1975
$ws.Range("J86").Value = 1980
$ws.Range("K86").Value = 1975
$ws.Range("L86").Value = 1980
$ws.Range("M86").Value = -852
$ws.Range("N86").Value = -4226
$ws.Range("H89").Value = 1977.7778
$ws.Range("I89").Value = 1975
$ws.Range("J89").Value = 1980
$ws.Range("K89").Value = 9875
$ws.Range("L89").Value = 9900
$ws.Range("M89").Value = -4259
$ws.Range("N89").Value = -21132
$ws.Range("H134").Value = 27893.928
$ws.Range("I134").Value = 4250.775
$ws.Range("J134").Value = 500757
$ws.Range("K134").Value = 12752.325
$ws.Range("L134").Value = 1502271
$ws.Range("M134").Value = -10217.325
$ws.Range("N134").Value = -1507341
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H16").Value = 2405.0356
$ws.Range("I16").Value = 2367.5908
$ws.Range("J16").Value = 2542.3333
$ws.Range("K16").Value = 2367.5908
$ws.Range("L16").Value = 2542.3333
$ws.Range("M16").Value = -2080.5908
$ws.Range("N16").Value = -3116.3333
$ws.Range("H31").Value = 16133269
$ws.Range("I31").Value = 1456.25
$ws.Range("K31").Value = 1456.25
$ws.Range("M31").Value = -1161.25
$ws.Range("H34").Value = 16133269
$ws.Range("I34").Value = 1456.25
$ws.Range("K34").Value = 1456.25
$ws.Range("M34").Value = -1254.25
$ws.Range("H94").Value = 3470.8572
$ws.Range("I94").Value = 1882.375
$ws.Range("J94").Value = 4448.385
$ws.Range("K94").Value = 1882.375
$ws.Range("L94").Value = 4448.385
$ws.Range("M94").Value = -1431.375
$ws.Range("N94").Value = -5350.385
$ws.Range("H113").Value = 2405.0356
$ws.Range("I113").Value = 2367.5908
$ws.Range("J113").Value = 2542.3333
$ws.Range("K113").Value = 2367.5908
$ws.Range("L113").Value = 2542.3333
$ws.Range("M113").Value = -197.5907999999999
$ws.Range("N113").Value = -6882.3333
$ws.Range("H122").Value = 4235.875
$ws.Range("I122").Value = 3012.4
$ws.Range("K122").Value = 9037.200000000001
$ws.Range("M122").Value = -6587.200000000001
$ws.Range("H132").Value = 5264686.5
$ws.Range("I132").Value = 6251034.5
$ws.Range("J132").Value = 4162.6665
$ws.Range("K132").Value = 18753103.5
$ws.Range("L132").Value = 12487.9995
$ws.Range("M132").Value = -18750573.5
$ws.Range("N132").Value = -17547.9995
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H12").Value = 3448399.5
$ws.Range("J12").Value = 117.75
$ws.Range("L12").Value = 353.25
$ws.Range("N12").Value = -699.25
$ws.Range("H80").Value = 3266.6667
$ws.Range("I80").Value = 3400
$ws.Range("K80").Value = 10200
$ws.Range("M80").Value = -9264
$ws.Range("H83").Value = 3266.6667
$ws.Range("I83").Value = 3400
$ws.Range("K83").Value = 30600
$ws.Range("M83").Value = -25920
$ws.Range("H122").Value = 6878.0557
$ws.Range("I122").Value = 1116.1818
$ws.Range("J122").Value = 15932.429
$ws.Range("K122").Value = 10045.6362
$ws.Range("L122").Value = 143391.861
$ws.Range("M122").Value = -7595.636200000001
$ws.Range("N122").Value = -148291.861
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H113").Value = 83342620
$ws.Range("I113").Value = 142872140
$ws.Range("K113").Value = 142872140
$ws.Range("M113").Value = -142869970
$ws.Range("H122").Value = 48744910
$ws.Range("I122").Value = 62637004
$ws.Range("J122").Value = 27275316
$ws.Range("K122").Value = 187911012
$ws.Range("L122").Value = 81825948
$ws.Range("M122").Value = -187908562
$ws.Range("N122").Value = -81830848
$ws.Range("H126").Value = 11999.211
$ws.Range("I126").Value = 14266
$ws.Range("J126").Value = 3498.75
$ws.Range("K126").Value = 42798
$ws.Range("L126").Value = 10496.25
$ws.Range("M126").Value = -40328
$ws.Range("N126").Value = -15436.25
$ws.Range("H132").Value = 7938661
$ws.Range("I132").Value = 9261049
$ws.Range("K132").Value = 27783147
$ws.Range("M132").Value = -27780617
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H7").Value = 2771.182
$ws.Range("I7").Value = 2047.8
$ws.Range("K7").Value = 2047.8
$ws.Range("M7").Value = -1935.8
$ws.Range("H40").Value = 2258
$ws.Range("I40").Value = 2134.3333
$ws.Range("K40").Value = 2134.3333
$ws.Range("M40").Value = -1998.3333
$ws.Range("H126").Value = 2771.182
$ws.Range("I126").Value = 2047.8
$ws.Range("K126").Value = 6143.4
$ws.Range("M126").Value = -3673.4
$ws.Range("H132").Value = 6669.1763
$ws.Range("I132").Value = 6891.483
$ws.Range("J132").Value = 5379.8
$ws.Range("K132").Value = 20674.449
$ws.Range("L132").Value = 16139.4
$ws.Range("M132").Value = -18144.449
$ws.Range("N132").Value = -21199.4
$ws.Range("H136").Value = 11434.458
$ws.Range("I136").Value = 9268.9375
$ws.Range("J136").Value = 15765.5
$ws.Range("K136").Value = 27806.8125
$ws.Range("L136").Value = 47296.5
$ws.Range("M136").Value = -25256.8125
$ws.Range("N136").Value = -52396.5
$ws.Range("H137").Value = 53333.332
$ws.Range("J137").Value = 55000
$ws.Range("L137").Value = 55000
$ws.Range("N137").Value = -65200
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H141").Value = 105815.71
$ws.Range("I141").Value = 220000
$ws.Range("J141").Value = 86785
$ws.Range("K141").Value = 220000
$ws.Range("L141").Value = 86785
$ws.Range("M141").Value = -214820
$ws.Range("N141").Value = -97145
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H122").Value = 1354.1538
$ws.Range("I122").Value = 1176
$ws.Range("J122").Value = 1433.3334
$ws.Range("K122").Value = 3528
$ws.Range("L122").Value = 4300.0002
$ws.Range("M122").Value = -1078
$ws.Range("N122").Value = -9200.0002
$ws.Range("H137").Value = 45715
$ws.Range("J137").Value = 45715
$ws.Range("L137").Value = 45715
$ws.Range("N137").Value = -55915
$ws.Range("H138").Value = 59095.668
$ws.Range("J138").Value = 59095.668
$ws.Range("L138").Value = 59095.668
$ws.Range("N138").Value = -69375.66800000001
$ws.Range("H139").Value = 59069
$ws.Range("J139").Value = 59069
$ws.Range("L139").Value = 59069
$ws.Range("N139").Value = -69349
$ws.Range("H141").Value = 78683.57000000001
$ws.Range("J141").Value = 78683.57000000001
$ws.Range("L141").Value = 78683.57000000001
$ws.Range("N141").Value = -89043.57000000001
